$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 (ID 9): R1, R12, R15, R19 -> R1, R15 ; Qty 4 -> 2
$ws.Range("B10").Value = "R1, R15"
$ws.Range("D10").Value = 2

# Row 11 (ID 10): R10 -> R10, R43 ; Qty 1 -> 2
$ws.Range("B11").Value = "R10, R43"
$ws.Range("D11").Value = 2

# Row 12 (ID 11): R11, R13, R14, R17, R18 -> R11, R14, R17 ; Qty 5 -> 3
$ws.Range("B12").Value = "R11, R14, R17"
$ws.Range("D12").Value = 3

# Row 16 (ID 15): R26, R27, R30, R33 -> R26, R27, R30, R33, R41, R42 ; Qty 4 -> 6
$ws.Range("B16").Value = "R26, R27, R30, R33, R41, R42"
$ws.Range("D16").Value = 6

# Row 20 (ID 19): R36, R37 -> R36, R37, R44 ; Qty 2 -> 3
$ws.Range("B20").Value = "R36, R37, R44"
$ws.Range("D20").Value = 3

# Row 22 (ID 21): R41, R42, R6, R7, R8, R9 -> R6, R7, R8, R9 ; Qty 6 -> 4
$ws.Range("B22").Value = "R6, R7, R8, R9"
$ws.Range("D22").Value = 4

# Row 23 (ID 22, U1): SOIC-8 -> SOIC-14 ; TL972IDR -> TL974IDR (Value + Manufacturer Part columns)
$ws.Range("C23").Value = "SOIC-14"
$ws.Range("E23").Value = "TL974IDR"
$ws.Range("F23").Value = "TL974IDR"

# Move active cell selection to F24 (cosmetic view-state change)
$ws.Range("F24").Select()
